# Update Titan_Profits market-price derived columns (H..N) across sheets
# per scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 592549.5
$ws.Range("J64").Value = 7300
$ws.Range("L64").Value = 7300
$ws.Range("N64").Value = -7796
$ws.Range("H67").Value = 592549.5
$ws.Range("J67").Value = 7300
$ws.Range("L67").Value = 7300
$ws.Range("N67").Value = -9016
$ws.Range("H129").Value = 1081.2174
$ws.Range("I129").Value = 495.25
$ws.Range("J129").Value = 1137.0238
$ws.Range("K129").Value = 1485.75
$ws.Range("L129").Value = 3411.0714
$ws.Range("M129").Value = 3514.25
$ws.Range("N129").Value = -13411.0714
$ws.Range("H137").Value = 20000954
$ws.Range("I137").Value = 24390946
$ws.Range("J137").Value = 2095.7778
$ws.Range("K137").Value = 73172838
$ws.Range("L137").Value = 6287.3334
$ws.Range("M137").Value = -73170288
$ws.Range("N137").Value = -11387.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2076.1042
$ws.Range("I61").Value = 1563.8292
$ws.Range("J61").Value = 5076.5713
$ws.Range("K61").Value = 1563.8292
$ws.Range("L61").Value = 5076.5713
$ws.Range("M61").Value = -1351.8292
$ws.Range("N61").Value = -5500.5713
$ws.Range("H63").Value = 18335.75
$ws.Range("I63").Value = 22563.334
$ws.Range("J63").Value = 5653
$ws.Range("K63").Value = 22563.334
$ws.Range("L63").Value = 5653
$ws.Range("M63").Value = -21877.334
$ws.Range("N63").Value = -7025
$ws.Range("H66").Value = 18335.75
$ws.Range("I66").Value = 22563.334
$ws.Range("J66").Value = 5653
$ws.Range("K66").Value = 112816.67
$ws.Range("L66").Value = 28265
$ws.Range("M66").Value = -109384.67
$ws.Range("N66").Value = -35129
$ws.Range("H132").Value = 3053
$ws.Range("I132").Value = 2480.8635
$ws.Range("J132").Value = 4851.143
$ws.Range("K132").Value = 7442.5905
$ws.Range("L132").Value = 14553.429
$ws.Range("M132").Value = -4912.5905
$ws.Range("N132").Value = -19613.429
$ws.Range("H136").Value = 2076.1042
$ws.Range("I136").Value = 1563.8292
$ws.Range("J136").Value = 5076.5713
$ws.Range("K136").Value = 4691.487599999999
$ws.Range("L136").Value = 15229.7139
$ws.Range("M136").Value = -2141.487599999999
$ws.Range("N136").Value = -20329.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087.6316
$ws.Range("I31").Value = 1236.3928
$ws.Range("J31").Value = 4471.1
$ws.Range("K31").Value = 1236.3928
$ws.Range("L31").Value = 4471.1
$ws.Range("M31").Value = -941.3928000000001
$ws.Range("N31").Value = -5061.1
$ws.Range("H34").Value = 2087.6316
$ws.Range("I34").Value = 1236.3928
$ws.Range("J34").Value = 4471.1
$ws.Range("K34").Value = 1236.3928
$ws.Range("L34").Value = 4471.1
$ws.Range("M34").Value = -1034.3928
$ws.Range("N34").Value = -4875.1
$ws.Range("H62").Value = 18982.691
$ws.Range("I62").Value = 19731.25
$ws.Range("K62").Value = 19731.25
$ws.Range("M62").Value = -19107.25
$ws.Range("H65").Value = 18982.691
$ws.Range("I65").Value = 19731.25
$ws.Range("K65").Value = 98656.25
$ws.Range("M65").Value = -95536.25
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 2404.6924
$ws.Range("I132").Value = 1938.5758
$ws.Range("K132").Value = 5815.7274
$ws.Range("M132").Value = -3285.7274
$ws.Range("H134").Value = 1968.5968
$ws.Range("I134").Value = 1257.1296
$ws.Range("J134").Value = 6771
$ws.Range("K134").Value = 3771.3888
$ws.Range("L134").Value = 20313
$ws.Range("M134").Value = -1236.3888
$ws.Range("N134").Value = -25383

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2719.6365
$ws.Range("I64").Value = 818.6667
$ws.Range("J64").Value = 5000.8
$ws.Range("K64").Value = 2456.0001
$ws.Range("L64").Value = 15002.4
$ws.Range("M64").Value = -2186.0001
$ws.Range("N64").Value = -15542.4
$ws.Range("H67").Value = 2719.6365
$ws.Range("I67").Value = 818.6667
$ws.Range("J67").Value = 5000.8
$ws.Range("K67").Value = 2456.0001
$ws.Range("L67").Value = 15002.4
$ws.Range("M67").Value = -1520.0001
$ws.Range("N67").Value = -16874.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 22000
$ws.Range("J62").Value = 22000
$ws.Range("L62").Value = 22000
$ws.Range("N62").Value = -23372
$ws.Range("H65").Value = 22000
$ws.Range("J65").Value = 22000
$ws.Range("L65").Value = 66000
$ws.Range("N65").Value = -72864
$ws.Range("H80").Value = 5771.4287
$ws.Range("I80").Value = 5984.615
$ws.Range("K80").Value = 5984.615
$ws.Range("M80").Value = -4986.615
$ws.Range("H83").Value = 5771.4287
$ws.Range("I83").Value = 5984.615
$ws.Range("K83").Value = 29923.075
$ws.Range("M83").Value = -24931.075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1899.7858
$ws.Range("I68").Value = 1399.6666
$ws.Range("J68").Value = 2800
$ws.Range("K68").Value = 1399.6666
$ws.Range("L68").Value = 2800
$ws.Range("M68").Value = -650.6666
$ws.Range("N68").Value = -4298
$ws.Range("H71").Value = 1899.7858
$ws.Range("I71").Value = 1399.6666
$ws.Range("J71").Value = 2800
$ws.Range("K71").Value = 6998.333000000001
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = -3254.333000000001
$ws.Range("N71").Value = -21488
$ws.Range("H82").Value = 2400
$ws.Range("I82").Value = 2200
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 2200
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1839
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 2400
$ws.Range("I85").Value = 2200
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 2200
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -952
$ws.Range("N85").Value = -4996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9617833
$ws.Range("I132").Value = 13891061
$ws.Range("K132").Value = 41673183
$ws.Range("M132").Value = -41670653
